$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Recorded By" - move "System" from front to back of the list
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, System"

# Row 3: "Recorded By" - reorder names
$ws.Range("G3").Value = "hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 9: "Recorded By" - reorder names
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Row 10: Average Attendance % statistic updated; keep the cell text-formatted
# (same as original) so it is not auto-converted to a numeric percentage.
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "27.3%"

# Row 15: add a recorder name, update attendance count and coverage %
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("H15").Value = "59/251"

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "27.3%"
